$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete row 6 (shrinks used range from AH6 to AH5)
$ws.Rows.Item(6).Delete()

# New data values for rows 2-5 (columns A..AH)
$data = @(
    @(45160.50694444445, 4.284, 5.265, 0.768, 10.407, 9.715, 5.079, 10.361, 6.858, 3.662, 6.216, 7.246, 4.306, 0.973, 3.629, 6.538, 2.613, 0.117, 0.053, 60.452, 12.195, 3.56, 7.364, 6.727, 0.825, 5.701, 3.6, 5.538, 2.243, 6.732, 0.372, 8.518000000000001, 2.17, 5.204),
    @(45160.51388888889, 4.962, 4.725, 0.514, 11.511, 10.038, 4.711, 15.63, 7.059, 3.446, 5.413, 6.216, 4.88, 1.176, 3.981, 6.564, 3.315, 0.063, 0.022, 60.97, 12.603, 3.761, 8.180999999999999, 5.556, 0.669, 8.266999999999999, 3.651, 4.289, 3.08, 5.985, 0.237, 14.194, 2.318, 5.206),
    @(45160.52083333334, 18.655, 14.685, 0.861, 41.126, 34.343, 15.209, 54.119, 23.363, 10.661, 16.088, 17.588, 17.472, 4.636, 14.66, 21.568, 12.227, 0.168, 0.438, 221.336, 41.959, 13.576, 28.203, 15.641, 2.01, 27.378, 12.247, 11.445, 12.193, 18.093, 0.175, 49.073, 7.891, 17.354),
    @(45160.52777777778, 8.69, 7.04, 0.47, 19.35, 16.25, 7.23, 30.91, 11.14, 5.15, 7.77, 8.59, 8.23, 2.15, 6.85, 10.3, 5.77, 0.08, 0.16, 101.25, 20.11, 6.35, 13.47, 7.63, 0.96, 14.94, 5.82, 5.62, 5.67, 8.720000000000001, 0.14, 28.32, 3.75, 8.24)
)

for ($r = 0; $r -lt $data.Length; $r++) {
    $rowValues = $data[$r]
    $excelRow = $r + 2
    for ($c = 0; $c -lt $rowValues.Length; $c++) {
        $ws.Cells.Item($excelRow, $c + 1).Value = $rowValues[$c]
    }
}

# Column width adjustments.
# The ColumnWidth COM property adds a fixed padding offset (~0.8333333333333334,
# i.e. 5/6 character) relative to the raw OOXML <col width="..."> value, so we
# compensate by subtracting that offset to land exactly on the target width.
$widthOffset = 0.8333333333333334
$ws.Columns.Item(7).ColumnWidth = 8 - $widthOffset   # G: 7 -> 8
$ws.Columns.Item(9).ColumnWidth = 8 - $widthOffset   # I: 7 -> 8
$ws.Columns.Item(10).ColumnWidth = 8 - $widthOffset  # J: 7 -> 8
$ws.Columns.Item(15).ColumnWidth = 7 - $widthOffset  # O: 8 -> 7
$ws.Columns.Item(27).ColumnWidth = 8 - $widthOffset  # AA: 7 -> 8
$ws.Columns.Item(28).ColumnWidth = 8 - $widthOffset  # AB: 7 -> 8
